$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column A text (lowercased chemical names, new header label) ---
$ws.Range("A1").Value = "chemical name"
$ws.Range("A2").Value = "polystyrene"
$ws.Range("A3").Value = "poly(methylmethcrylate)"
$ws.Range("A4").Value = "poly(2-vinylpyridine)"
$ws.Range("A5").Value = "poly(ethylmethcrylate)"

# --- Bold the header row, with A1 slightly larger than B1:D1 ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 12
$ws.Range("B1:D1").Font.Bold = $true

# --- Move the active selection to H5 ---
[void]$ws.Range("H5").Select()
